$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2:B4 hold numeric-looking job numbers that must stay TEXT (shared-string)
# cells, matching the source data's existing storage type. A leading
# apostrophe forces text entry so the value isn't re-interpreted as a number.
$ws.Range("B2").Value = "'32297400"
$ws.Range("B3").Value = "'32297401"
$ws.Range("B4").Value = "'32297402"
